# Apply the two text-structure edits from the commit to the "Introduction"
# slide's body placeholder:
#   1. " Separate your Concerns" (single run, leading space) becomes two
#      runs: "Separate " + "your Concerns" (leading space dropped).
#   2. "Extensible : Ability to create new " + "custom directives" (two
#      runs, trailing endParaRPr) becomes a single run "Extensible :
#      Ability to create new custom directives" with no endParaRPr.

$p = $ppt.ActivePresentation

# Locate the slide/shape containing the target paragraphs instead of
# assuming fixed indexes.
$sh = $null
for ($k = 1; $k -le $p.Slides.Count; $k++) {
    $slide = $p.Slides.Item($k)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $shapeTxt = $shape.TextFrame.TextRange.Text
            if ($shapeTxt -match "Extensible" -and $shapeTxt -match "Separate your Concerns") {
                $sh = $shape
            }
        }
    }
}

$tf = $sh.TextFrame
$tr = $tf.TextRange

# Locate the target paragraphs by content instead of hard-coded indexes.
$paraCount = $tr.Paragraphs().Count
$separateIdx   = 0
$extensibleIdx = 0
$nextIdx       = 0
for ($i = 1; $i -le $paraCount; $i++) {
    $txt = $tr.Paragraphs($i, 1).Text
    if ($txt -match "Separate your Concerns") {
        $separateIdx = $i
    }
    if ($txt -match "^Extensible") {
        $extensibleIdx = $i
        $nextIdx = $i + 1
    }
}

# --- Edit 1: split " Separate your Concerns" into "Separate " + "your Concerns" ---
$paraSep = $tr.Paragraphs($separateIdx, 1)
$paraSep.Text = "Separate your Concerns"
$firstPart = $tr.Characters($paraSep.Start, 9)
$firstPart.Text = "Separate "

# --- Edit 2: merge "Extensible : Ability to create new " + "custom directives"
#             into a single run and drop the paragraph's endParaRPr ---
# The endParaRPr is carried by the paragraph mark; temporarily delete the
# following paragraph (remembering its text) so the merge/rewrite below
# doesn't leave a stray endParaRPr behind, then reinsert it unchanged.
$paraNext = $tr.Paragraphs($nextIdx, 1)
# Paragraph-scoped .Text includes the trailing paragraph-mark (CR); strip
# it since InsertAfter below re-adds the break explicitly.
$nextText = $paraNext.Text.TrimEnd([char]13)
$paraNext.Delete()

$paraExt = $tr.Paragraphs($extensibleIdx, 1)
$paraExt.Text = [guid]::NewGuid().ToString()
$paraExt2 = $tr.Paragraphs($extensibleIdx, 1)
$paraExt2.Text = "Extensible : Ability to create new custom directives"
$paraExt3 = $tr.Paragraphs($extensibleIdx, 1)
$null = $paraExt3.InsertAfter([char]13 + $nextText)
